$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the "Ver no Jupiter ..." paragraph (Jekyll-site footer line).
$jupiterIndex = -1
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Ver no Jupiter*") {
        $jupiterIndex = $i
        break
    }
}

if ($jupiterIndex -gt 0) {
    # The blank separator paragraph right before it, and the copyright
    # paragraph right after it, are removed together with it.
    $startPara = $paras.Item($jupiterIndex - 1)
    $endPara   = $paras.Item($jupiterIndex + 1)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
